$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column R (values were 0, now filled with actuals)
$ws.Range("R2").Value = 8100.92
$ws.Range("R3").Value = 3454
$ws.Range("R4").Value = 2358
$ws.Range("R5").Value = 3704.15
$ws.Range("R6").Value = 17617.07

# Update column AG (totals) to reflect the new sums
$ws.Range("AG2").Value = 151943.76
$ws.Range("AG3").Value = 75376.00999999999
$ws.Range("AG4").Value = 51739.9
$ws.Range("AG5").Value = 48272.02
$ws.Range("AG6").Value = 327331.69

$wb.Save()
